# The workbook holds a small "create local user" table on Sheet1 with
# columns UserName/FullName/Description/Password/CanUserChangeThePassword/
# Operation/Servers (A:G) across rows 2-4.
#
# The edit replaces every "Servers" value (column G, rows 2-4) -- which used
# to reference the old 192.168.176.x addresses -- with the single new
# address 192.168.4.92, and leaves the final selection on G4 (the last cell
# touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "192.168.4.92"
$ws.Range("G3").Value = "192.168.4.92"
$ws.Range("G4").Value = "192.168.4.92"

$ws.Range("G4").Select()
